$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")
$ws.Rows("1:2").Delete() | Out-Null
$ws.Rows("1:1").Select() | Out-Null
